$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Rows 37 / 38 had their match data (columns F:V) swapped.
#    Columns A:E (index/pais/torneio/temporada/data_partida) stay put.
# ------------------------------------------------------------------
$row37 = $ws.Range("F37:V37").Value2
$row38 = $ws.Range("F38:V38").Value2
$ws.Range("F37:V37").Value2 = $row38
$ws.Range("F38:V38").Value2 = $row37

# ------------------------------------------------------------------
# 2) Rows 39 / 40 had their match data (columns F:V) swapped too.
# ------------------------------------------------------------------
$row39 = $ws.Range("F39:V39").Value2
$row40 = $ws.Range("F40:V40").Value2
$ws.Range("F39:V39").Value2 = $row40
$ws.Range("F40:V40").Value2 = $row39

# ------------------------------------------------------------------
# 3) A brand-new match row (45) was appended, mirroring the layout and
#    formatting of the last existing row (44).
# ------------------------------------------------------------------
$ws.Range("A44:V44").Copy()
$ws.Range("A45:V45").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A45").Value2 = 44
$ws.Range("B45").Value2 = "united-arab-emirates"
$ws.Range("C45").Value2 = "uae-league"
$ws.Range("D45").Value2 = "2023-2024"
$ws.Range("E45").Value2 = 45232.6875
$ws.Range("F45").Value2 = "Al Wasl"
$ws.Range("G45").Value2 = 3
$ws.Range("H45").Value2 = "Al Ain"
$ws.Range("I45").Value2 = 1
$ws.Range("J45").Value2 = 2.43
$ws.Range("K45").Value2 = "28/10/2023 15:13"
$ws.Range("L45").Value2 = 2.64
$ws.Range("M45").Value2 = "02/11/2023 16:27"
$ws.Range("N45").Value2 = 3.83
$ws.Range("O45").Value2 = "28/10/2023 15:13"
$ws.Range("P45").Value2 = 3.83
$ws.Range("Q45").Value2 = "02/11/2023 16:27"
$ws.Range("R45").Value2 = 2.51
$ws.Range("S45").Value2 = "28/10/2023 15:13"
$ws.Range("T45").Value2 = 2.47
$ws.Range("U45").Value2 = "02/11/2023 16:27"
$ws.Range("V45").Value2 = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-wasl-al-ain/zaUm54eT/"

Write-Output "Edit complete"
